$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new score values for row 27 (weeport_final2 entry)
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 1
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 9
$ws.Range("I27").Value = 13
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 2

# Update the frozen-pane view: scroll the unfrozen pane down so row 24 is at
# its top, and move the active selection to L27.
$win = $excel.ActiveWindow
$win.ScrollRow = 24
$win.ScrollColumn = 1
$ws.Range("L27").Select()
